# Big stimulus update: replace the "face" image category with "book",
# and expand the abbreviated correct_ans codes in column L to their full words.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# ---------------------------------------------------------------------------
# 1) Stimulus filename swap: "face//face_NN.jpg" -> "book//book_NN.jpg"
#    Occurs in columns A (promptFile), B (correctFile), C (dist_01File) and
#    D (dist_02File).
# ---------------------------------------------------------------------------
$fileCols = 1, 2, 3, 4   # A, B, C, D

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $fileCols) {
        $cell = $ws.Cells.Item($r, $c)
        $cur = $cell.Value2
        if ($cur -and $cur.ToString().StartsWith("face//face_")) {
            $newVal = $cur.ToString().Replace("face//face_", "book//book_")
            $cell.Value = $newVal
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Expand abbreviated correct_ans codes in column L (12):
#       b -> center
#       y -> left
#       r -> right
# ---------------------------------------------------------------------------
$ansCol = 12   # L

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $ansCol)
    $cur = $cell.Value2
    if ($cur -eq "b") {
        $cell.Value = "center"
    } elseif ($cur -eq "y") {
        $cell.Value = "left"
    } elseif ($cur -eq "r") {
        $cell.Value = "right"
    }
}
